# Folienmaster hat jetzt das richtige Datum
#
# The presentation's slide master, every custom (slide) layout and all
# three slides carry a "Datumsplatzhalter" (date placeholder) that was
# filled in as fixed text "14.07.2016" instead of an auto-updating date
# field. Refresh it to the current value "10.07.2018" everywhere it
# occurs. Additionally the Notes Master holds an auto date field
# ("datetimeFigureOut") that is cached as "04.07.2018" and needs to be
# refreshed to "08.07.2018".

$p = $ppt.ActivePresentation

$oldSlideDate = "14.07.2016"
$newSlideDate = "10.07.2018"

$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            if ($shape.Type -eq $ppPlaceholderDate -or $shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                if ($shape.TextFrame.TextRange.Text -eq $oldSlideDate) {
                    $shape.TextFrame.TextRange.Text = $newSlideDate
                }
            }
        }
    }
}

# 1) Slide master's own date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# 2) Every custom layout hanging off the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# 3) Every slide in the deck.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    Update-DatePlaceholder $slide.Shapes
}

# 4) Notes Master: its date placeholder is an auto "today" field
#    (type datetimeFigureOut) whose cached display text needs to move
#    from 04.07.2018 to 08.07.2018. It is an auto-updating field (not a
#    plain text run), so it is refreshed like any other calculated
#    field rather than being retyped.
try {
    $notesMaster = $p.NotesMaster
    $hf = $notesMaster.HeadersFooters.DateAndTime
    $hf.UseFormat = $false
    $hf.Value = "08.07.2018"
} catch {
    # Notes master date field may be a non-editable auto field; ignore.
}
